$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("first_eval")

# The evaluation window has advanced by one quarter: every existing row's
# B:G values shift down into the next row (Q6 row keeps its Q-label but now
# shows what used to be Q7's numbers, etc.), and the top row (Q6) is
# refreshed with the newly computed error metrics for the latest quarter.
# Column G (N) simply continues the existing count-down sequence by one.

$newRow2 = @(-0.02314597604078636, 0.3579920056255013, 0.1782699060034266, 0.4222202103209018, 0.4363822494547141, 15)

for ($r = 11; $r -ge 3; $r--) {
    for ($c = 2; $c -le 7; $c++) {
        $ws.Cells.Item($r, $c).Value = $ws.Cells.Item($r - 1, $c).Value2
    }
}

for ($c = 2; $c -le 7; $c++) {
    $ws.Cells.Item(2, $c).Value = $newRow2[$c - 2]
}
